$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only column Z needs to be cleared
$onlyZRows = @(57,58,59,60,61,62,63,65,66,67,68,69,70,71,72,73,74,75,76,77,78,80)
foreach ($r in $onlyZRows) {
    $ws.Range("Z$r").ClearContents()
}

# Rows where the entire data range (B:AA) needs to be cleared
$fullRows = @(64,79)
foreach ($r in $fullRows) {
    $ws.Range("B$r`:AA$r").ClearContents()
}
